# ExcelColumnDateValidator now supports Locale: add two more example rows
# to the "dates" sheet showing German-locale formatted date / date-time
# strings (to be used as test data for the locale-aware validator).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dates")

# Make sure the new string cells are stored/rendered as text (same as the
# existing text example rows 4-5) rather than being re-interpreted as a
# date/number.
$ws.Range("A6:B7").NumberFormat = "@"

$ws.Range("A6").Value = "21. November 2020"
$ws.Range("B6").Value = "21. November 2020 01:02"

$ws.Range("A7").Value = "21. November 2020"
$ws.Range("B7").Value = "21. November 2020 1:2:17"

# Column A now holds text long enough to need a wider column.
$ws.Columns.Item(1).ColumnWidth = 17.33203125

# Selection moves to the last entered cell, as left by the editor.
$ws.Range("A7").Select()
